# training_diary.xlsx - add support for preset & sticky values (closes #32)
#
# Adds a new "value" column (J) to the "survey" sheet that lets an item
# specify a preset/sticky value:
#   - most items just get the literal keyword "sticky" (re-use the value
#     that was entered last time)
#   - the squats item instead gets a formula-ish preset head(squats,1)
#     (pull the most recent squats entry)
# Also refreshes the mood smiley icons to the outline/3x variants and
# bumps the mood row's height down now that the extra column fits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# --- new column header -----------------------------------------------
$ws.Range("J1").Value = "value"

# --- per-item preset / sticky values -----------------------------------
$ws.Range("J3").Value = "sticky"          # pushups
$ws.Range("J4").Value = "sticky"          # situps
$ws.Range("J5").Value = "sticky"          # pullups
$ws.Range("J6").Value = "head(squats,1)"  # squats -> pull last entry
$ws.Range("J7").Value = "sticky"          # jogging_km
$ws.Range("J8").Value = "sticky"          # jogging_min

# J4 gets its own left-aligned / wrap-text style
$ws.Range("J4").HorizontalAlignment = -4131
$ws.Range("J4").WrapText = $true

# --- refresh the mood icons (outline + 3x variants) ---------------------
$ws.Range("G10").Value = '<i class="fa fa-smile-o fa-3x"></i>'
$ws.Range("H10").Value = '<i class="fa fa-frown-o fa-3x"></i>'

# mood row shrinks a bit now that J fits without extra wrapping
$ws.Rows.Item(10).RowHeight = 45

# --- restore selection / view -------------------------------------------
$ws.Activate()
$ws.Range("H16").Select() | Out-Null
